# Update the "cryptos" price/volume table with the latest scraped values.
# Values that could be mistaken for numbers (e.g. "214.85") are entered with a
# leading apostrophe so Excel keeps them as literal text, matching the
# original inline-string cell contents exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.829.47"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "1.633.82"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "'214.85"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").Value = "'0.5030"
$ws.Range("E6").Value = "  -2.28%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "'0.2573"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "'0.06406"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").Value = "'19.64"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").Value = "'0.07663"
$ws.Range("E11").Value = "  -1.82%  "
$ws.Range("D12").Value = "1.636.85"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").Value = "'4.238"
$ws.Range("E13").Value = "  -1.52%  "
$ws.Range("D14").Value = "1.857.69"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "'0.5458"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").Value = "0.0₅7925"
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").Value = "'63.51"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").Value = "25.825.82"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "'202.98"
$ws.Range("E20").Value = "  -4.08%  "
$ws.Range("D21").Value = "'4.313"
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("D22").Value = "'9.936"
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("D23").Value = "'5.961"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D26").Value = "'141.00"
$ws.Range("E26").Value = "  -2.27%  "
$ws.Range("D27").Value = "'0.1144"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").Value = "'15.71"
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("D29").Value = "'6.690"
$ws.Range("E29").Value = "  -4.24%  "
$ws.Range("D30").Value = "'1.239"
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("D31").Value = "'0.04978"
$ws.Range("E31").Value = "  -4.47%  "
$ws.Range("D32").Value = "'3.275"
$ws.Range("E32").Value = "  -2.71%  "
$ws.Range("D33").Value = "'3.182"
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D34").Value = "'1.533"
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("D36").Value = "1.173.46"
$ws.Range("D37").Value = "'0.8923"
$ws.Range("E37").Value = "  -4.33%  "
$ws.Range("D38").Value = "'2.615"
$ws.Range("E38").Value = "  -5.32%  "
$ws.Range("D39").Value = "'0.5577"
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").Value = "'5.647"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").Value = "'99.36"
$ws.Range("D45").Value = "'0.8020"
$ws.Range("E45").Value = "  -5.24%  "
$ws.Range("D46").Value = "1.770.24"
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "'0.4510"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Value = "'54.79"
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("E51").Value = "  -0.75%  "
